$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A218").NumberFormat = "@"
$ws.Range("A218").Value = "2026-02-06"
$ws.Range("B218").Value = "09:58:10"
$ws.Range("C218").Value = "09:00"
$ws.Range("D218").Value = "Bathroom"
$ws.Range("E218").Value = "No Motion"
$ws.Range("F218").Value = "Inactive"

$ws.Range("A219").NumberFormat = "@"
$ws.Range("A219").Value = "2026-02-06"
$ws.Range("B219").Value = "09:58:12"
$ws.Range("C219").Value = "09:00"
$ws.Range("D219").Value = "Bathroom"
$ws.Range("E219").Value = "No Motion"
$ws.Range("F219").Value = "Inactive"

$ws.Range("A220").NumberFormat = "@"
$ws.Range("A220").Value = "2026-02-06"
$ws.Range("B220").Value = "09:58:15"
$ws.Range("C220").Value = "09:00"
$ws.Range("D220").Value = "Bathroom"
$ws.Range("E220").Value = "No Motion"
$ws.Range("F220").Value = "Inactive"

$ws.Range("A221").NumberFormat = "@"
$ws.Range("A221").Value = "2026-02-06"
$ws.Range("B221").Value = "09:58:20"
$ws.Range("C221").Value = "09:00"
$ws.Range("D221").Value = "Bathroom"
$ws.Range("E221").Value = "No Motion"
$ws.Range("F221").Value = "Inactive"

$ws.Range("A222").NumberFormat = "@"
$ws.Range("A222").Value = "2026-02-06"
$ws.Range("B222").Value = "09:58:25"
$ws.Range("C222").Value = "09:00"
$ws.Range("D222").Value = "Bathroom"
$ws.Range("E222").Value = "No Motion"
$ws.Range("F222").Value = "Inactive"

$ws.Range("A223").NumberFormat = "@"
$ws.Range("A223").Value = "2026-02-06"
$ws.Range("B223").Value = "09:58:30"
$ws.Range("C223").Value = "09:00"
$ws.Range("D223").Value = "Bathroom"
$ws.Range("E223").Value = "No Motion"
$ws.Range("F223").Value = "Inactive"

$ws.Range("A224").NumberFormat = "@"
$ws.Range("A224").Value = "2026-02-06"
$ws.Range("B224").Value = "09:58:35"
$ws.Range("C224").Value = "09:00"
$ws.Range("D224").Value = "Bathroom"
$ws.Range("E224").Value = "No Motion"
$ws.Range("F224").Value = "Inactive"

$ws.Range("A225").NumberFormat = "@"
$ws.Range("A225").Value = "2026-02-06"
$ws.Range("B225").Value = "09:58:40"
$ws.Range("C225").Value = "09:00"
$ws.Range("D225").Value = "Bathroom"
$ws.Range("E225").Value = "No Motion"
$ws.Range("F225").Value = "Inactive"

$ws.Range("A226").NumberFormat = "@"
$ws.Range("A226").Value = "2026-02-06"
$ws.Range("B226").Value = "09:58:45"
$ws.Range("C226").Value = "09:00"
$ws.Range("D226").Value = "Bathroom"
$ws.Range("E226").Value = "No Motion"
$ws.Range("F226").Value = "Inactive"

$ws.Range("A227").NumberFormat = "@"
$ws.Range("A227").Value = "2026-02-06"
$ws.Range("B227").Value = "09:58:50"
$ws.Range("C227").Value = "09:00"
$ws.Range("D227").Value = "Bathroom"
$ws.Range("E227").Value = "No Motion"
$ws.Range("F227").Value = "Inactive"

$ws.Range("A228").NumberFormat = "@"
$ws.Range("A228").Value = "2026-02-06"
$ws.Range("B228").Value = "09:58:55"
$ws.Range("C228").Value = "09:00"
$ws.Range("D228").Value = "Bathroom"
$ws.Range("E228").Value = "No Motion"
$ws.Range("F228").Value = "Inactive"

$ws.Range("A229").NumberFormat = "@"
$ws.Range("A229").Value = "2026-02-06"
$ws.Range("B229").Value = "09:59:00"
$ws.Range("C229").Value = "09:00"
$ws.Range("D229").Value = "Bathroom"
$ws.Range("E229").Value = "No Motion"
$ws.Range("F229").Value = "Inactive"

$ws.Range("A230").NumberFormat = "@"
$ws.Range("A230").Value = "2026-02-06"
$ws.Range("B230").Value = "09:59:06"
$ws.Range("C230").Value = "09:00"
$ws.Range("D230").Value = "Bathroom"
$ws.Range("E230").Value = "No Motion"
$ws.Range("F230").Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A133").NumberFormat = "@"
$ws.Range("A133").Value = "2026-02-06"
$ws.Range("B133").Value = "09:58:11"
$ws.Range("C133").Value = "09:00"
$ws.Range("D133").Value = "Bathroom"
$ws.Range("E133").NumberFormat = "@"
$ws.Range("E133").Value = "70.1%"
$ws.Range("F133").Value = "Active"

$ws.Range("A134").NumberFormat = "@"
$ws.Range("A134").Value = "2026-02-06"
$ws.Range("B134").Value = "09:58:13"
$ws.Range("C134").Value = "09:00"
$ws.Range("D134").Value = "Bathroom"
$ws.Range("E134").NumberFormat = "@"
$ws.Range("E134").Value = "69.0%"
$ws.Range("F134").Value = "Active"

$ws.Range("A135").NumberFormat = "@"
$ws.Range("A135").Value = "2026-02-06"
$ws.Range("B135").Value = "09:58:23"
$ws.Range("C135").Value = "09:00"
$ws.Range("D135").Value = "Bathroom"
$ws.Range("E135").NumberFormat = "@"
$ws.Range("E135").Value = "68.8%"
$ws.Range("F135").Value = "Active"

$ws.Range("A136").NumberFormat = "@"
$ws.Range("A136").Value = "2026-02-06"
$ws.Range("B136").Value = "09:58:28"
$ws.Range("C136").Value = "09:00"
$ws.Range("D136").Value = "Bathroom"
$ws.Range("E136").NumberFormat = "@"
$ws.Range("E136").Value = "69.8%"
$ws.Range("F136").Value = "Active"

$ws.Range("A137").NumberFormat = "@"
$ws.Range("A137").Value = "2026-02-06"
$ws.Range("B137").Value = "09:58:33"
$ws.Range("C137").Value = "09:00"
$ws.Range("D137").Value = "Bathroom"
$ws.Range("E137").NumberFormat = "@"
$ws.Range("E137").Value = "68.8%"
$ws.Range("F137").Value = "Active"

$ws.Range("A138").NumberFormat = "@"
$ws.Range("A138").Value = "2026-02-06"
$ws.Range("B138").Value = "09:58:43"
$ws.Range("C138").Value = "09:00"
$ws.Range("D138").Value = "Bathroom"
$ws.Range("E138").NumberFormat = "@"
$ws.Range("E138").Value = "68.8%"
$ws.Range("F138").Value = "Active"

$ws.Range("A139").NumberFormat = "@"
$ws.Range("A139").Value = "2026-02-06"
$ws.Range("B139").Value = "09:58:48"
$ws.Range("C139").Value = "09:00"
$ws.Range("D139").Value = "Bathroom"
$ws.Range("E139").NumberFormat = "@"
$ws.Range("E139").Value = "69.9%"
$ws.Range("F139").Value = "Active"

$ws.Range("A140").NumberFormat = "@"
$ws.Range("A140").Value = "2026-02-06"
$ws.Range("B140").Value = "09:58:58"
$ws.Range("C140").Value = "09:00"
$ws.Range("D140").Value = "Bathroom"
$ws.Range("E140").NumberFormat = "@"
$ws.Range("E140").Value = "69.8%"
$ws.Range("F140").Value = "Active"

$ws.Range("A141").NumberFormat = "@"
$ws.Range("A141").Value = "2026-02-06"
$ws.Range("B141").Value = "09:59:08"
$ws.Range("C141").Value = "09:00"
$ws.Range("D141").Value = "Bathroom"
$ws.Range("E141").NumberFormat = "@"
$ws.Range("E141").Value = "69.9%"
$ws.Range("F141").Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A132").NumberFormat = "@"
$ws.Range("A132").Value = "2026-02-06"
$ws.Range("B132").Value = "09:58:10"
$ws.Range("C132").Value = "09:00"
$ws.Range("D132").Value = "Bathroom"
$ws.Range("E132").Value = "27.9C"
$ws.Range("F132").Value = "Active"

$ws.Range("A133").NumberFormat = "@"
$ws.Range("A133").Value = "2026-02-06"
$ws.Range("B133").Value = "09:58:12"
$ws.Range("C133").Value = "09:00"
$ws.Range("D133").Value = "Bathroom"
$ws.Range("E133").Value = "27.9C"
$ws.Range("F133").Value = "Active"

$ws.Range("A134").NumberFormat = "@"
$ws.Range("A134").Value = "2026-02-06"
$ws.Range("B134").Value = "09:58:14"
$ws.Range("C134").Value = "09:00"
$ws.Range("D134").Value = "Bathroom"
$ws.Range("E134").Value = "27.8C"
$ws.Range("F134").Value = "Active"

$ws.Range("A135").NumberFormat = "@"
$ws.Range("A135").Value = "2026-02-06"
$ws.Range("B135").Value = "09:58:24"
$ws.Range("C135").Value = "09:00"
$ws.Range("D135").Value = "Bathroom"
$ws.Range("E135").Value = "27.9C"
$ws.Range("F135").Value = "Active"

$ws.Range("A136").NumberFormat = "@"
$ws.Range("A136").Value = "2026-02-06"
$ws.Range("B136").Value = "09:58:29"
$ws.Range("C136").Value = "09:00"
$ws.Range("D136").Value = "Bathroom"
$ws.Range("E136").Value = "27.9C"
$ws.Range("F136").Value = "Active"

$ws.Range("A137").NumberFormat = "@"
$ws.Range("A137").Value = "2026-02-06"
$ws.Range("B137").Value = "09:58:34"
$ws.Range("C137").Value = "09:00"
$ws.Range("D137").Value = "Bathroom"
$ws.Range("E137").Value = "27.8C"
$ws.Range("F137").Value = "Active"

$ws.Range("A138").NumberFormat = "@"
$ws.Range("A138").Value = "2026-02-06"
$ws.Range("B138").Value = "09:58:44"
$ws.Range("C138").Value = "09:00"
$ws.Range("D138").Value = "Bathroom"
$ws.Range("E138").Value = "27.8C"
$ws.Range("F138").Value = "Active"

$ws.Range("A139").NumberFormat = "@"
$ws.Range("A139").Value = "2026-02-06"
$ws.Range("B139").Value = "09:58:49"
$ws.Range("C139").Value = "09:00"
$ws.Range("D139").Value = "Bathroom"
$ws.Range("E139").Value = "27.8C"
$ws.Range("F139").Value = "Active"

$ws.Range("A140").NumberFormat = "@"
$ws.Range("A140").Value = "2026-02-06"
$ws.Range("B140").Value = "09:58:59"
$ws.Range("C140").Value = "09:00"
$ws.Range("D140").Value = "Bathroom"
$ws.Range("E140").Value = "27.8C"
$ws.Range("F140").Value = "Active"

$ws.Range("A141").NumberFormat = "@"
$ws.Range("A141").Value = "2026-02-06"
$ws.Range("B141").Value = "09:59:09"
$ws.Range("C141").Value = "09:00"
$ws.Range("D141").Value = "Bathroom"
$ws.Range("E141").Value = "27.8C"
$ws.Range("F141").Value = "Active"
